$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Rows 1-12: simple single-run text replacements
$table.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$table.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$table.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$table.Rows.Item(4).Cells.Item(1).Range.Text = "1807"
$table.Rows.Item(5).Cells.Item(1).Range.Text = "0.00001"
$table.Rows.Item(6).Cells.Item(1).Range.Text = "0.00073"
$table.Rows.Item(7).Cells.Item(1).Range.Text = "0.00011"
$table.Rows.Item(9).Cells.Item(1).Range.Text = "0.00015"
$table.Rows.Item(10).Cells.Item(1).Range.Text = "0.00016"
$table.Rows.Item(11).Cells.Item(1).Range.Text = "0.00017"
$table.Rows.Item(12).Cells.Item(1).Range.Text = "0.21602"

# Rows 44-46: collapse the multi-run tab-separated content down to a
# single value (these rows previously held a full tab-separated stats
# line; now they hold just the single summary number).
$table.Rows.Item(44).Cells.Item(1).Range.Text = "99.89"
$table.Rows.Item(45).Cells.Item(1).Range.Text = "0.22"
$table.Rows.Item(46).Cells.Item(1).Range.Text = "200"
